$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 0.07225593691223042
$ws.Range("J4").Value = 0.5349644508647857
$ws.Range("K4").Value = 0.7136139080027292
$ws.Range("L4").Value = 2.972965183273405
